$wb = $excel.ActiveWorkbook

# --- Rename existing headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins used on the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header styling (bold, bordered, centered) used on the other sheets
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45067.99999999999, 232, -88.95374758580186, 586.3986332774882),
    @(45081.99999999999, 229, -121.2431143027604, 572.2215619308997),
    @(45088.99999999999, 227, -115.6440075829759, 556.2279496782595),
    @(45095.99999999999, 225, -89.65735796664218, 580.7689286694443),
    @(45102.99999999999, 223, -135.968968837059, 535.5951357787567),
    @(45109.99999999999, 222, -140.8198938272831, 564.082258517399),
    @(45116.99999999999, 220, -117.6203288753938, 559.4574266713938),
    @(45123.99999999999, 218, -120.6546781992352, 547.9232538822794),
    @(45130.99999999999, 216, -113.2993014735265, 564.5956327519781),
    @(45144.99999999999, 212, -132.0451004741186, 544.3509718307134),
    @(45151.99999999999, 211, -144.472661488273, 547.7369280398793),
    @(45165.99999999999, 207, -106.9525462149501, 555.340666880013),
    @(45179.99999999999, 203, -115.7692149484944, 536.5810986093268),
    @(45186.99999999999, 202, -121.1468490585957, 545.8408219429367),
    @(45200.99999999999, 198, -119.1746683686142, 545.4356273500005),
    @(45207.99999999999, 196, -151.8343082641204, 542.1021462248227),
    @(45242.99999999999, 187, -143.8057940983222, 524.4749421940367),
    @(45249.99999999999, 185, -156.2175152864835, 514.3258888971727),
    @(45277.99999999999, 178, -169.7354629290046, 495.7109267495015),
    @(45305.99999999999, 171, -170.909688784987, 503.3438918974319),
    @(45312.99999999999, 169, -152.5292379729926, 504.7238684154186),
    @(45319.99999999999, 167, -162.6231851905745, 503.6828206650554),
    @(45333.99999999999, 163, -191.2217008059024, 497.4488993907822),
    @(45340.99999999999, 161, -172.6220497694241, 492.6341003416475),
    @(45466.99999999999, 129, -198.3494707290961, 471.2291136597877),
    @(45480.99999999999, 125, -203.8556790421866, 492.907496990975),
    @(45487.99999999999, 123, -194.1942585671445, 456.5910547270826),
    @(45501.99999999999, 119, -250.5999209227412, 448.4630977578096),
    @(45522.99999999999, 114, -236.5504672935728, 470.7995043904547),
    @(45529.99999999999, 112, -224.9673293424748, 469.7017998100837),
    @(45536.99999999999, 110, -233.8215146733762, 442.2252375079244),
    @(45550.99999999999, 107, -245.7097018722992, 435.0277587962874),
    @(45557.99999999999, 105, -218.8941271740536, 425.2312488504497),
    @(45564.99999999999, 103, -245.8794554193183, 416.310362106894),
    @(45578.99999999999, 99, -241.3867850299935, 442.0035864291083),
    @(45599.99999999999, 94, -239.5442918078082, 425.7101638092473),
    @(45606.99999999999, 92, -257.9191798329433, 415.084038883354),
    @(45613.99999999999, 90, -251.3265184723479, 416.6406902092184),
    @(45620.99999999999, 89, -266.151354314655, 445.5920191557597),
    @(45627.99999999999, 87, -235.2738410456996, 413.0288899739418),
    @(45634.99999999999, 85, -282.0628476219604, 422.2555410882775),
    @(45641.99999999999, 83, -220.4844129237707, 396.4469210384087),
    @(45648.99999999999, 81, -268.4046648917768, 397.0657187533897),
    @(45655.99999999999, 79, -308.1331336268926, 416.5245151797489)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match date styling (column A) used on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A45").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore original active sheet/selection state
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
